# Adapt column header formatting to respective input file names (#7)
# - rename "<Column>_old" headers to "<Column>_FV2310"
# - rename "<Column>_new" headers to "<Column>_FV2404"
# - turn the data range into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1, columns A:U) -----------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $name = $cell.Value()
    if ($name -like "*_old") {
        $cell.Value = ($name.Substring(0, $name.Length - 4) + "_FV2310")
    } elseif ($name -like "*_new") {
        $cell.Value = ($name.Substring(0, $name.Length - 4) + "_FV2404")
    }
}

# --- Convert the data range into a real Excel table ----------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U63"), $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row ------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
